# ResINSOR.xlsx result update:
# The dependent-value ("DepVal") recorded for the newInsOR_001 test case is
# updated from the previous ticket reference to the new one, reflecting the
# latest alert-message validation run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "EAOR21AP-0315"
